# Fix incorrect combining of dictionaries when fetching existing meetings
# and webinars from zoom: append "+iclrzoom" to the local part of each
# alternate host / panelist email address in column J.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("J2").Value = "Timnit Gebru <tgebru+iclrzoom@gmail.com>"
$ws.Range("J3").Value = "Shakir Mohamed <shakir.mohamed+iclrzoom@gmail.com>"
$ws.Range("J4").Value = "Kyunghyun Cho <kyunghyun.cho+iclrzoom@nyu.edu>"
$ws.Range("J5").Value = "Asja Fischer <asja.fischer+iclrzoom@gmail.com>"
$ws.Range("J6").Value = "Martha White <whitem+iclrzoom@ualberta.ca>"
$ws.Range("J7").Value = "Gabriel Synnaeve <gabriel.synnaeve+iclrzoom@gmail.com>"
$ws.Range("J8").Value = "Dawn Song <dawnsong+iclrzoom@gmail.com>"
$ws.Range("J9").Value = "Alexander Rush <sasha.rush+iclrzoom@gmail.com>"

# Column J needs to be widened to fit the longer addresses (37 -> ~51.66 chars).
$ws.Columns.Item(10).ColumnWidth = 50.83

# Move the active selection to J11 (matches the saved cursor position).
$ws.Range("J11").Select() | Out-Null
